$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.7
$ws.Range("I2").Value = 5.75
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.3
$ws.Range("X2").Value = 7
$ws.Range("AH2").Value = 13
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 8.5
$ws.Range("AW2").Value = 7
$ws.Range("AY2").Value = 41
$ws.Range("BA2").Value = 151

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 2.7
$ws.Range("V3").Value = 1.58
$ws.Range("AU3").Value = 9

# Row 4
$ws.Range("G4").Value = 1.66
$ws.Range("R4").Value = 1.62
$ws.Range("AT4").Value = 2.62

# Row 5
$ws.Range("G5").Value = 2.82
$ws.Range("I5").Value = 2.35
$ws.Range("M5").Value = 1.05
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.92
$ws.Range("R5").Value = 1.82

# Row 6
$ws.Range("G6").Value = 1.44
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 5.6
$ws.Range("J6").Value = 1.9
$ws.Range("K6").Value = 2.4
$ws.Range("L6").Value = 5.4
$ws.Range("P6").Value = 3.94
$ws.Range("R6").Value = 2.07
$ws.Range("W6").Value = 6.7
$ws.Range("X6").Value = 6.3
$ws.Range("Z6").Value = 8.5
$ws.Range("AA6").Value = 9.5
$ws.Range("AB6").Value = 18.5
$ws.Range("AC6").Value = 13.5
$ws.Range("AD6").Value = 7.4
$ws.Range("AE6").Value = 14
$ws.Range("AF6").Value = 55
$ws.Range("AH6").Value = 14.5
$ws.Range("AI6").Value = 29
$ws.Range("AJ6").Value = 15
$ws.Range("AK6").Value = 80
$ws.Range("AL6").Value = 40
$ws.Range("AM6").Value = 40
$ws.Range("AN6").Value = 3.35
$ws.Range("AO6").Value = 6.5
$ws.Range("AQ6").Value = 18
$ws.Range("AS6").Value = 175
$ws.Range("AT6").Value = 3.2
$ws.Range("AU6").Value = 7.7
$ws.Range("AW6").Value = 7.4
$ws.Range("AX6").Value = 32
$ws.Range("AY6").Value = 32
$ws.Range("AZ6").Value = 200
$ws.Range("BA6").Value = 200
$ws.Range("BB6").Value = 400

# Row 8
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.25
$ws.Range("BC8").Value = 151
$ws.Range("BD8").Value = 151

# Row 9
$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.65
$ws.Range("R9").Value = 2.2

# Row 10
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.75
$ws.Range("J10").Value = 2.75
$ws.Range("K10").Value = 2.05
$ws.Range("L10").Value = 4.33
$ws.Range("M10").Value = 1.07
$ws.Range("O10").Value = 1.36
$ws.Range("Q10").Value = 2.15
$ws.Range("R10").Value = 1.67
$ws.Range("U10").Value = 1.87
$ws.Range("V10").Value = 1.77
$ws.Range("X10").Value = 9
$ws.Range("Z10").Value = 17
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 13
$ws.Range("AN10").Value = 4
$ws.Range("AR10").Value = 67
$ws.Range("AV10").Value = 51
$ws.Range("AW10").Value = 5.5
$ws.Range("AX10").Value = 21
$ws.Range("AZ10").Value = 67

# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63
$ws.Range("V11").Value = 1.63

# Row 12
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 3.25
$ws.Range("J12").Value = 3
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("V12").Value = 1.63
$ws.Range("AC12").Value = 7
$ws.Range("AK12").Value = 41
$ws.Range("AL12").Value = 34
$ws.Range("AQ12").Value = 41

# Row 13
$ws.Range("M13").Value = 1.05
$ws.Range("O13").Value = 1.25
$ws.Range("V13").Value = 1.58
